{"js": "// Update the date header and the 100 arithmetic-problem cells in the\n// practice-sheet table. Cells are matched by their current (old) text so\n// the script fails loudly instead of silently mis-writing a cell if the\n// document doesn't look like what we expect.\n\nconst DATE_OLD = \"2025-09-11 Thursday\";\nconst DATE_NEW = \"2025-09-12 Friday\";\n\n// [oldText, newText] for every table cell, in row-major (reading) order.\nconst CELL_REPLACEMENTS = [\n  [\"63-45=\", \"13+69=\"],\n  [\"30-15=\", \"25-12=\"],\n  [\"97-65=\", \"66-4=\"],\n  [\"94-29=\", \"94-48=\"],\n  [\"21+42=\", \"0+93=\"],\n  [\"86-70=\", \"9+5=\"],\n  [\"57-26=\", \"76-41=\"],\n  [\"10+29=\", \"11+83=\"],\n  [\"87-14=\", \"66-26=\"],\n  [\"82-14=\", \"93-31=\"],\n  [\"70-9=\", \"30+62=\"],\n  [\"97-61=\", \"65-56=\"],\n  [\"0+32=\", \"88-62=\"],\n  [\"48+18=\", \"52+33=\"],\n  [\"54-30=\", \"12+12=\"],\n  [\"39+55=\", \"48-4=\"],\n  [\"70-39=\", \"59+13=\"],\n  [\"17-7=\", \"66-47=\"],\n  [\"10+59=\", \"87+5=\"],\n  [\"6+14=\", \"51+37=\"],\n  [\"30+55=\", \"15+18=\"],\n  [\"10+1=\", \"69-7=\"],\n  [\"44-40=\", \"61+38=\"],\n  [\"36+44=\", \"86-59=\"],\n  [\"1+25=\", \"28+34=\"],\n  [\"30-11=\", \"9+46=\"],\n  [\"36+7=\", \"75-71=\"],\n  [\"99-39=\", \"39-25=\"],\n  [\"98-58=\", \"26+12=\"],\n  [\"80-4=\", \"13-1=\"],\n  [\"85-41=\", \"16+30=\"],\n  [\"17+58=\", \"32-9=\"],\n  [\"42+55=\", \"49-13=\"],\n  [\"91-45=\", \"14+13=\"],\n  [\"10+69=\", \"96-88=\"],\n  [\"2+96=\", \"51-9=\"],\n  [\"94-47=\", \"22+64=\"],\n  [\"66-65=\", \"30+59=\"],\n  [\"3+92=\", \"51-48=\"],\n  [\"9+16=\", \"86-32=\"],\n  [\"54+4=\", \"73-21=\"],\n  [\"61+4=\", \"47-37=\"],\n  [\"97-27=\", \"26+19=\"],\n  [\"47-19=\", \"72+22=\"],\n  [\"30+1=\", \"36+28=\"],\n  [\"50+15=\", \"63-3=\"],\n  [\"5+27=\", \"1+46=\"],\n  [\"26+72=\", \"42+17=\"],\n  [\"6+2=\", \"28-18=\"],\n  [\"4+84=\", \"59+22=\"],\n  [\"7+44=\", \"32-16=\"],\n  [\"22+77=\", \"8+74=\"],\n  [\"96-45=\", \"73-67=\"],\n  [\"96-18=\", \"47-3=\"],\n  [\"76-15=\", \"24+1=\"],\n  [\"1+48=\", \"76-5=\"],\n  [\"23-18=\", \"85-63=\"],\n  [\"83-23=\", \"24+18=\"],\n  [\"66-10=\", \"8+49=\"],\n  [\"10-6=\", \"60+36=\"],\n  [\"58+18=\", \"59+23=\"],\n  [\"85-1=\", \"67+3=\"],\n  [\"8+73=\", \"32+42=\"],\n  [\"65+31=\", \"99-7=\"],\n  [\"2+15=\", \"53-25=\"],\n  [\"28+66=\", \"14+3=\"],\n  [\"78-23=\", \"18+75=\"],\n  [\"12+22=\", \"62-59=\"],\n  [\"92-6=\", \"12+70=\"],\n  [\"68-34=\", \"20+13=\"],\n  [\"37+27=\", \"22+8=\"],\n  [\"86-68=\", \"17+34=\"],\n  [\"36+52=\", \"12+13=\"],\n  [\"57+32=\", \"46+34=\"],\n  [\"41+32=\", \"44+44=\"],\n  [\"53+17=\", \"35+46=\"],\n  [\"86-40=\", \"46-18=\"],\n  [\"30-9=\", \"0+68=\"],\n  [\"54-15=\", \"10+31=\"],\n  [\"78-46=\", \"11-0=\"],\n  [\"23+64=\", \"20+38=\"],\n  [\"31-13=\", \"19+73=\"],\n  [\"83-30=\", \"19+35=\"],\n  [\"95-53=\", \"72-12=\"],\n  [\"3+17=\", \"66+15=\"],\n  [\"61-9=\", \"22+62=\"],\n  [\"50+16=\", \"49+27=\"],\n  [\"20+69=\", \"60-29=\"],\n  [\"99-25=\", \"4+66=\"],\n  [\"15+30=\", \"15+71=\"],\n  [\"65-7=\", \"52+33=\"],\n  [\"48+22=\", \"26+53=\"],\n  [\"93-29=\", \"5+67=\"],\n  [\"60-53=\", \"41+49=\"],\n  [\"78-62=\", \"44+51=\"],\n  [\"86-82=\", \"97-37=\"],\n  [\"83-79=\", \"58-31=\"],\n  [\"12+29=\", \"55-32=\"],\n  [\"41+55=\", \"90-29=\"],\n  [\"85-29=\", \"22+54=\"]\n];\n\n// --- 1. Update the date paragraph (first paragraph in the body) -----------\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst dateParagraph = paragraphs.items[0];\ndateParagraph.load(\"text\");\nawait context.sync();\n\nif (dateParagraph.text !== DATE_OLD) {\n  throw new Error(\n    `Unexpected date paragraph text: expected \"${DATE_OLD}\", found \"${dateParagraph.text}\"`\n  );\n}\ndateParagraph.insertText(DATE_NEW, \"Replace\");\nawait context.sync();\n\n// --- 2. Update every cell in the first table -------------------------------\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.rows.load(\"items\");\nawait context.sync();\n\nconst rows = table.rows.items;\nfor (const row of rows) {\n  row.cells.load(\"items\");\n}\nawait context.sync();\n\n// Flatten all cells in row-major order and load their current values.\nconst cells = [];\nfor (const row of rows) {\n  for (const cell of row.cells.items) {\n    cell.load(\"value\");\n    cells.push(cell);\n  }\n}\nawait context.sync();\n\nif (cells.length !== CELL_REPLACEMENTS.length) {\n  throw new Error(\n    `Expected ${CELL_REPLACEMENTS.length} table cells, found ${cells.length}`\n  );\n}\n\nfor (let i = 0; i < cells.length; i++) {\n  const [oldText, newText] = CELL_REPLACEMENTS[i];\n  const cell = cells[i];\n  if (cell.value !== oldText) {\n    throw new Error(\n      `Cell ${i} text mismatch: expected \"${oldText}\", found \"${cell.value}\"`\n    );\n  }\n  cell.value = newText;\n}\nawait context.sync();\n", "ps1": "# Update the date header and the 100 arithmetic-problem cells in the\n# practice-sheet table. Cells are matched by their current (old) text so\n# the script fails loudly instead of silently mis-writing a cell if the\n# document doesn't look like what we expect.\n\n$d = $word.ActiveDocument\n\n$dateOld = '2025-09-11 Thursday'\n$dateNew = '2025-09-12 Friday'\n\n# [oldText, newText] for every table cell, in row-major (reading) order.\n$cellReplacements = @(\n    @('63-45=', '13+69='),\n    @('30-15=', '25-12='),\n    @('97-65=', '66-4='),\n    @('94-29=', '94-48='),\n    @('21+42=', '0+93='),\n    @('86-70=', '9+5='),\n    @('57-26=', '76-41='),\n    @('10+29=', '11+83='),\n    @('87-14=', '66-26='),\n    @('82-14=', '93-31='),\n    @('70-9=', '30+62='),\n    @('97-61=', '65-56='),\n    @('0+32=', '88-62='),\n    @('48+18=', '52+33='),\n    @('54-30=', '12+12='),\n    @('39+55=', '48-4='),\n    @('70-39=', '59+13='),\n    @('17-7=', '66-47='),\n    @('10+59=', '87+5='),\n    @('6+14=', '51+37='),\n    @('30+55=', '15+18='),\n    @('10+1=', '69-7='),\n    @('44-40=', '61+38='),\n    @('36+44=', '86-59='),\n    @('1+25=', '28+34='),\n    @('30-11=', '9+46='),\n    @('36+7=', '75-71='),\n    @('99-39=', '39-25='),\n    @('98-58=', '26+12='),\n    @('80-4=', '13-1='),\n    @('85-41=', '16+30='),\n    @('17+58=', '32-9='),\n    @('42+55=', '49-13='),\n    @('91-45=', '14+13='),\n    @('10+69=', '96-88='),\n    @('2+96=', '51-9='),\n    @('94-47=', '22+64='),\n    @('66-65=', '30+59='),\n    @('3+92=', '51-48='),\n    @('9+16=', '86-32='),\n    @('54+4=', '73-21='),\n    @('61+4=', '47-37='),\n    @('97-27=', '26+19='),\n    @('47-19=', '72+22='),\n    @('30+1=', '36+28='),\n    @('50+15=', '63-3='),\n    @('5+27=', '1+46='),\n    @('26+72=', '42+17='),\n    @('6+2=', '28-18='),\n    @('4+84=', '59+22='),\n    @('7+44=', '32-16='),\n    @('22+77=', '8+74='),\n    @('96-45=', '73-67='),\n    @('96-18=', '47-3='),\n    @('76-15=', '24+1='),\n    @('1+48=', '76-5='),\n    @('23-18=', '85-63='),\n    @('83-23=', '24+18='),\n    @('66-10=', '8+49='),\n    @('10-6=', '60+36='),\n    @('58+18=', '59+23='),\n    @('85-1=', '67+3='),\n    @('8+73=', '32+42='),\n    @('65+31=', '99-7='),\n    @('2+15=', '53-25='),\n    @('28+66=', '14+3='),\n    @('78-23=', '18+75='),\n    @('12+22=', '62-59='),\n    @('92-6=', '12+70='),\n    @('68-34=', '20+13='),\n    @('37+27=', '22+8='),\n    @('86-68=', '17+34='),\n    @('36+52=', '12+13='),\n    @('57+32=', '46+34='),\n    @('41+32=', '44+44='),\n    @('53+17=', '35+46='),\n    @('86-40=', '46-18='),\n    @('30-9=', '0+68='),\n    @('54-15=', '10+31='),\n    @('78-46=', '11-0='),\n    @('23+64=', '20+38='),\n    @('31-13=', '19+73='),\n    @('83-30=', '19+35='),\n    @('95-53=', '72-12='),\n    @('3+17=', '66+15='),\n    @('61-9=', '22+62='),\n    @('50+16=', '49+27='),\n    @('20+69=', '60-29='),\n    @('99-25=', '4+66='),\n    @('15+30=', '15+71='),\n    @('65-7=', '52+33='),\n    @('48+22=', '26+53='),\n    @('93-29=', '5+67='),\n    @('60-53=', '41+49='),\n    @('78-62=', '44+51='),\n    @('86-82=', '97-37='),\n    @('83-79=', '58-31='),\n    @('12+29=', '55-32='),\n    @('41+55=', '90-29='),\n    @('85-29=', '22+54=')\n)\n\n# --- 1. Update the date paragraph (first paragraph in the body) -----------\n$dateParagraph = $d.Paragraphs.Item(1)\n$currentDateText = $dateParagraph.Range.Text.TrimEnd([char]13, [char]7)\nif ($currentDateText -ne $dateOld) {\n    throw \"Unexpected date paragraph text: expected '$dateOld', found '$currentDateText'\"\n}\n$dateParagraph.Range.Text = $dateNew\n\n# --- 2. Update every cell in the first table -------------------------------\n$tbl = $d.Tables.Item(1)\n\n$numCols = 5\nfor ($i = 0; $i -lt $cellReplacements.Count; $i++) {\n    $oldText = $cellReplacements[$i][0]\n    $newText = $cellReplacements[$i][1]\n\n    $row = [int][math]::Floor($i / $numCols) + 1\n    $col = ($i % $numCols) + 1\n\n    $cell = $tbl.Cell($row, $col)\n    $currentText = $cell.Range.Text.TrimEnd([char]13, [char]7)\n    if ($currentText -ne $oldText) {\n        throw \"Cell ($row,$col) text mismatch: expected '$oldText', found '$currentText'\"\n    }\n    $cell.Range.Text = $newText\n}\n"}
